$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values (row 1): B1=0 (existing), C1=1, D1=2 - all styled like the existing B1 (style index 1 / bold+border+center)
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 2

# Copy B1's formatting (bold font, thin border, centered) onto the new header cells
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

# Data values for columns C (col index 3) and D (col index 4), rows 2-25
$cValues = @(
    0.09180947066839784,
    0.09106907171139463,
    0.08588627901237217,
    0.08681177770862618,
    0.08662667796937538,
    0.08662667796937538,
    0.08366508214136255,
    0.08403528161986415,
    0.08588627901237217,
    0.08088858605260052,
    1.057104610861331,
    7.595197600678162,
    13.42287779125041,
    15.74495402015173,
    16.60696350584271,
    17.25,
    16.67711630701877,
    14.63768737995343,
    7.235178607835352,
    0.665063363128132,
    0.1155022372925005,
    0.10846844720097,
    0.1012495573701888,
    0.09884326075992833
)

$dValues = @(
    39.27027815818423,
    39.76939137462676,
    39.21276070314078,
    37.75910917540354,
    36.68639292369171,
    36.68639292369171,
    36.48649619468828,
    36.17910782839054,
    37.18707765529937,
    37.89080214624071,
    37.58184226477785,
    37.23736614058326,
    36.74171025750398,
    36.00561255416115,
    36.72505219675369,
    37.65916081090182,
    37.56895584042386,
    37.98697887434614,
    37.90494578272681,
    38.3776575443953,
    39.44408773544664,
    40.59443683631547,
    41.40376714635296,
    42
)

for ($i = 0; $i -lt $cValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}
